$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "444/444"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "IR801997"
$ws.Range("D2").Value = "NOUBAIL MOHAMMED"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 22500
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "--"
$ws.Range("O2").Value = 22500

# Row 3
$ws.Range("A3").Value = "444/444"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "B219321"
$ws.Range("D3").Value = "JEMAA HORMI"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "--"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "--"
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 22500
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 22500

# Row 4
$ws.Range("A4").Value = "444/444"
$ws.Range("B4").Value = "Direction régionale"
$ws.Range("C4").Value = "IR801997"
$ws.Range("D4").Value = "NOUBAIL MOHAMMED"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "--"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "--"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 3000

# Row 5
$ws.Range("A5").Value = "444/444"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "B219321"
$ws.Range("D5").Value = "JEMAA HORMI"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "--"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "--"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "--"
$ws.Range("O5").Value = 3000

# Row 6
$ws.Range("A6").Value = "444/444"
$ws.Range("B6").Value = "Direction régionale"
$ws.Range("C6").Value = "IR801997"
$ws.Range("D6").Value = "NOUBAIL MOHAMMED"
$ws.Range("E6").Value = "non"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 1500

# Row 7
$ws.Range("A7").Value = "444/444"
$ws.Range("B7").Value = "Direction régionale"
$ws.Range("C7").Value = "B219321"
$ws.Range("D7").Value = "JEMAA HORMI"
$ws.Range("E7").Value = "non"
$ws.Range("F7").Value = "mensuelle"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 1500

# Row 8
$ws.Range("A8").Value = "000/CCCC/AV1"
$ws.Range("B8").Value = "Direction régionale"
$ws.Range("C8").Value = "BK646476"
$ws.Range("D8").Value = "DOUNIA LAMKADDAM"
$ws.Range("E8").Value = "non"
$ws.Range("F8").Value = "mensuelle"
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = "--"
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "--"
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 32000
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = "--"
$ws.Range("O8").Value = 32000

# Row 9
$ws.Range("A9").Value = "000/CCCC/AV1"
$ws.Range("B9").Value = "Direction régionale"
$ws.Range("C9").Value = "BK646476"
$ws.Range("D9").Value = "DOUNIA LAMKADDAM"
$ws.Range("E9").Value = "non"
$ws.Range("F9").Value = "mensuelle"
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = "--"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "--"
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = "--"
$ws.Range("O9").Value = 4000

# Row 10
$ws.Range("A10").Value = "000/CCCC/AV1"
$ws.Range("B10").Value = "Direction régionale"
$ws.Range("C10").Value = "BK646476"
$ws.Range("D10").Value = "DOUNIA LAMKADDAM"
$ws.Range("E10").Value = "non"
$ws.Range("F10").Value = "mensuelle"
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = "--"
$ws.Range("O10").Value = 2000

# Row 11
$ws.Range("A11").Value = "555/RRR/AV10"
$ws.Range("B11").Value = "Direction régionale"
$ws.Range("C11").Value = "B171710"
$ws.Range("D11").Value = "NADIA BADRANE"
$ws.Range("E11").Value = "non"
$ws.Range("F11").Value = "mensuelle"
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = "--"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "--"
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = "--"
$ws.Range("O11").Value = 1800

# Row 12
$ws.Range("A12").Value = "555/RRR/AV10"
$ws.Range("B12").Value = "Direction régionale"
$ws.Range("C12").Value = "IB43905"
$ws.Range("D12").Value = "NHILA BELGACEM"
$ws.Range("E12").Value = "non"
$ws.Range("F12").Value = "mensuelle"
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = "--"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "--"
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = "--"
$ws.Range("O12").Value = 1800

# Row 13
$ws.Range("A13").Value = "555/RRR/AV10"
$ws.Range("B13").Value = "Direction régionale"
$ws.Range("C13").Value = "B171710"
$ws.Range("D13").Value = "NADIA BADRANE"
$ws.Range("E13").Value = "non"
$ws.Range("F13").Value = "mensuelle"
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 8500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 850
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = "--"
$ws.Range("O13").Value = 7650

# Row 14
$ws.Range("A14").Value = "555/RRR/AV10"
$ws.Range("B14").Value = "Direction régionale"
$ws.Range("C14").Value = "IB43905"
$ws.Range("D14").Value = "NHILA BELGACEM"
$ws.Range("E14").Value = "non"
$ws.Range("F14").Value = "mensuelle"
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 8500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 850
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = "--"
$ws.Range("O14").Value = 7650

# Row 15
$ws.Range("A15").Value = "555/RRR/AV10"
$ws.Range("B15").Value = "Direction régionale"
$ws.Range("C15").Value = "B171710"
$ws.Range("D15").Value = "NADIA BADRANE"
$ws.Range("E15").Value = "non"
$ws.Range("F15").Value = "mensuelle"
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = "--"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "--"
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = "--"
$ws.Range("O15").Value = 1800

# Row 16
$ws.Range("A16").Value = "555/RRR/AV10"
$ws.Range("B16").Value = "Direction régionale"
$ws.Range("C16").Value = "IB43905"
$ws.Range("D16").Value = "NHILA BELGACEM"
$ws.Range("E16").Value = "non"
$ws.Range("F16").Value = "mensuelle"
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = "--"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "--"
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = "--"
$ws.Range("O16").Value = 1800

# Row 17
$ws.Range("A17").Value = " "
$ws.Range("B17").Value = " "
$ws.Range("C17").Value = " "
$ws.Range("D17").Value = " "
$ws.Range("E17").Value = " "
$ws.Range("F17").Value = " "
$ws.Range("G17").Value = " "
$ws.Range("H17").Value = 22000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1700
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 95000
$ws.Range("M17").Value = 800
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 114500
